$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.863243333333333
$ws.Range("H2").Value = 5.589729999999999
$ws.Range("I2").Value = 0.6067417803684044
$ws.Range("J2").Value = 0.6067417803684044
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1003616666666667
$ws.Range("N2").Value = 0.301085
$ws.Range("O2").Value = 0.0420091445250749
$ws.Range("P2").Value = 0.04200914452507489
$ws.Range("Q2").Value = 0.1869982063388889
$ws.Range("R2").Value = 1.68298385705
$ws.Range("S2").Value = 0.02548870314089755
$ws.Range("T2").Value = 0.02548870314089755

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.863243333333333
$ws.Range("H3").Value = 5.589729999999999
$ws.Range("I3").Value = 0.6067417803684044
$ws.Range("J3").Value = 0.6067417803684044
$ws.Range("O3").Value = 0.5857450611332571
$ws.Range("P3").Value = 0.585745061133257
$ws.Range("Q3").Value = 2.607367444447778
$ws.Range("R3").Value = 23.46630700003
$ws.Range("S3").Value = 0.3553960012339923
$ws.Range("T3").Value = 0.3553960012339922

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.863243333333333
$ws.Range("H4").Value = 5.589729999999999
$ws.Range("I4").Value = 0.6067417803684044
$ws.Range("J4").Value = 0.6067417803684044
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.7630573333333333
$ws.Range("N4").Value = 2.289172
$ws.Range("O4").Value = 0.3193986993399032
$ws.Range("P4").Value = 0.3193986993399031
$ws.Range("Q4").Value = 1.421761489284444
$ws.Range("R4").Value = 12.79585340356
$ws.Range("S4").Value = 0.1937925354848456
$ws.Range("T4").Value = 0.1937925354848455

$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 1.863243333333333
$ws.Range("H5").Value = 5.589729999999999
$ws.Range("I5").Value = 0.6067417803684044
$ws.Range("J5").Value = 0.6067417803684044
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.126254
$ws.Range("N5").Value = 0.378762
$ws.Range("O5").Value = 0.05284709500176502
$ws.Range("P5").Value = 0.05284709500176501
$ws.Range("Q5").Value = 0.2352419238066666
$ws.Range("R5").Value = 2.11717731426
$ws.Range("S5").Value = 0.03206454050866911
$ws.Range("T5").Value = 0.03206454050866911

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.9841983333333334
$ws.Range("H6").Value = 2.952595
$ws.Range("I6").Value = 0.3204918210730839
$ws.Range("J6").Value = 0.3204918210730839
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1003616666666667
$ws.Range("N6").Value = 0.301085
$ws.Range("O6").Value = 0.0420091445250749
$ws.Range("P6").Value = 0.04200914452507489
$ws.Range("Q6").Value = 0.0987757850638889
$ws.Range("R6").Value = 0.888982065575
$ws.Range("S6").Value = 0.01346358723056363
$ws.Range("T6").Value = 0.01346358723056362

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.9841983333333334
$ws.Range("H7").Value = 2.952595
$ws.Range("I7").Value = 0.3204918210730839
$ws.Range("J7").Value = 0.3204918210730839
$ws.Range("O7").Value = 0.5857450611332571
$ws.Range("P7").Value = 0.585745061133257
$ws.Range("Q7").Value = 1.377257949782778
$ws.Range("R7").Value = 12.395321548045
$ws.Range("S7").Value = 0.1877265013271625
$ws.Range("T7").Value = 0.1877265013271624

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.9841983333333334
$ws.Range("H8").Value = 2.952595
$ws.Range("I8").Value = 0.3204918210730839
$ws.Range("J8").Value = 0.3204918210730839
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.7630573333333333
$ws.Range("N8").Value = 2.289172
$ws.Range("O8").Value = 0.3193986993399032
$ws.Range("P8").Value = 0.3193986993399031
$ws.Range("Q8").Value = 0.7509997557044444
$ws.Range("R8").Value = 6.75899780134
$ws.Range("S8").Value = 0.10236467079982
$ws.Range("T8").Value = 0.1023646707998199

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.9841983333333334
$ws.Range("H9").Value = 2.952595
$ws.Range("I9").Value = 0.3204918210730839
$ws.Range("J9").Value = 0.3204918210730839
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.126254
$ws.Range("N9").Value = 0.378762
$ws.Range("O9").Value = 0.05284709500176502
$ws.Range("P9").Value = 0.05284709500176501
$ws.Range("Q9").Value = 0.1242589763766667
$ws.Range("R9").Value = 1.11833078739
$ws.Range("S9").Value = 0.01693706171553794
$ws.Range("T9").Value = 0.01693706171553794

$ws.Range("G10").Value = 0.2234583333333333
$ws.Range("H10").Value = 0.670375
$ws.Range("I10").Value = 0.07276639855851162
$ws.Range("J10").Value = 0.07276639855851162
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.1003616666666667
$ws.Range("N10").Value = 0.301085
$ws.Range("O10").Value = 0.0420091445250749
$ws.Range("P10").Value = 0.04200914452507489
$ws.Range("Q10").Value = 0.02242665076388889
$ws.Range("R10").Value = 0.201839856875
$ws.Range("S10").Value = 0.003056854153613716
$ws.Range("T10").Value = 0.003056854153613716

$ws.Range("G11").Value = 0.2234583333333333
$ws.Range("H11").Value = 0.670375
$ws.Range("I11").Value = 0.07276639855851162
$ws.Range("J11").Value = 0.07276639855851162
$ws.Range("O11").Value = 0.5857450611332571
$ws.Range("P11").Value = 0.585745061133257
$ws.Range("Q11").Value = 0.3127009624027778
$ws.Range("R11").Value = 2.814308661625
$ws.Range("S11").Value = 0.04262255857210234
$ws.Range("T11").Value = 0.04262255857210233

$ws.Range("G12").Value = 0.2234583333333333
$ws.Range("H12").Value = 0.670375
$ws.Range("I12").Value = 0.07276639855851162
$ws.Range("J12").Value = 0.07276639855851162
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.7630573333333333
$ws.Range("N12").Value = 2.289172
$ws.Range("O12").Value = 0.3193986993399032
$ws.Range("P12").Value = 0.3193986993399031
$ws.Range("Q12").Value = 0.1705115199444444
$ws.Range("R12").Value = 1.5346036795
$ws.Range("S12").Value = 0.02324149305523762
$ws.Range("T12").Value = 0.02324149305523761

$ws.Range("G13").Value = 0.2234583333333333
$ws.Range("H13").Value = 0.670375
$ws.Range("I13").Value = 0.07276639855851162
$ws.Range("J13").Value = 0.07276639855851162
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.126254
$ws.Range("N13").Value = 0.378762
$ws.Range("O13").Value = 0.05284709500176502
$ws.Range("P13").Value = 0.05284709500176501
$ws.Range("Q13").Value = 0.02821250841666667
$ws.Range("R13").Value = 0.25391257575
$ws.Range("S13").Value = 0.003845492777557961
$ws.Range("T13").Value = 0.00384549277755796
